$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text so numeric-looking strings
# (e.g. "601.14") are not auto-converted to numbers by COM value assignment,
# matching the source file where these cells are stored as inline strings.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.061.54"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "3.519.99"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "601.14"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "184.03"
$ws.Range("E6").Value = "  +5.96%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.140"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").Value = "7.13"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").Value = "0.434"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "4.124.85"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "32.40"
$ws.Range("E13").Value = "  +12.46%  "
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "68.021.20"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "0.0000182"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "3.503.18"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "6.39"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "14.75"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").Value = "396.73"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "8.07"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "73.82"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "0.545"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "5.73"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "10.39"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").Value = "0.180"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "6.29"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "23.99"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "7.43"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").Value = "163.83"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").Value = "0.877"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").Value = "7.16"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").Value = "4.76"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "27.82"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "26.77"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "2.68"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0735"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.829.58"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "42.45"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").Value = "345.66"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "1.09"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "33.71"
$ws.Range("E51").Value = "  -0.32%  "

# Restore default (no explicit number format) styling on the touched range
# so the cells end up with the same style index as before the edit.
$dataRange.ClearFormats()

